$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(1519661118, 1, " Singapore_300mm_NAND ", " - ", " 250 ", " 1000 ", " Clean Room Manufacturing, R&D ", " Mon_Feb_26_11:05:18_EST_2018"),
    @(1519661136, 1, " Singapore_300mm_NAND ", " - ", " 250 ", " 1000 ", " Clean Room Manufacturing, R&D ", " Mon_Feb_26_11:05:36_EST_2018"),
    @(1519664173, 1, " Singapore_300mm_NAND ", " - ", " 250 ", " 100 ",  " Clean Room Manufacturing, R&D ", " Mon_Feb_26_11:56:13_EST_2018"),
    @(1519664511, 1, " Singapore_Hqs ",         " - ", " 250 ", " 999 ",  " Clean Room Manufacturing, R&D ", " Mon_Feb_26_12:01:51_EST_2018"),
    @(1519664549, 1, " Singapore_300mm_NAND ", " - ", " 250 ", " 1000 ", " Clean Room Manufacturing, R&D ", " Mon_Feb_26_12:02:29_EST_2018"),
    @(1519664962, 1, " Singapore_300mm_NAND ", " - ", " 250 ", " 10000 ", " Clean Room Manufacturing, R&D ", " Mon_Feb_26_12:09:22_EST_2018"),
    @(1519667496, 1, " Singapore_300mm_NAND ", " - ", " 250 ", " 9999 ", " Clean Room Manufacturing, R&D ", " Mon_Feb_26_12:51:36_EST_2018"),
    @(1519667630, 1, " Singapore_300mm_NAND ", " - ", " 250 ", " 9999 ", " Clean Room Manufacturing, R&D ", " Mon_Feb_26_12:53:50_EST_2018")
)

$startRow = 109
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]

    # Columns E and F look numeric ("250", "1000", ...); force text storage
    # so the leading/trailing padding spaces survive, then restore the
    # default "Normal" style so no stray style definition is introduced.
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 5).Style = "Normal"

    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 6).Style = "Normal"

    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
}
